$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.88"
$ws.Range("E2").Value = "'5.39%"
$ws.Range("E3").Value = "'12.99%"
$ws.Range("D4").Value = "'5.160"
$ws.Range("E4").Value = "'4.72%"
$ws.Range("D5").Value = "'0.07809"
$ws.Range("E5").Value = "'6.29%"
$ws.Range("D6").Value = "'2.398"
$ws.Range("E6").Value = "'8.60%"
$ws.Range("D7").Value = "'8.053"
$ws.Range("E7").Value = "'4.36%"
$ws.Range("D8").Value = "'3.965"
$ws.Range("E8").Value = "'6.19%"
$ws.Range("D9").Value = "'0.9315"
$ws.Range("E9").Value = "'3.05%"
$ws.Range("D10").Value = "'0.09988"
$ws.Range("E10").Value = "'8.65%"
$ws.Range("D11").Value = "'0.1843"
$ws.Range("E11").Value = "'9.70%"
$ws.Range("D12").Value = "'0.08663"
$ws.Range("E12").Value = "'4.58%"
$ws.Range("D13").Value = "'0.03316"
$ws.Range("E13").Value = "'6.30%"
$ws.Range("D14").Value = "'0.09908"
$ws.Range("D15").Value = "'0.001494"
$ws.Range("E15").Value = "'-0.21%"
$ws.Range("D16").Value = "'0.005741"
$ws.Range("E16").Value = "'-0.23%"
$ws.Range("D17").Value = "'3.469"
$ws.Range("E17").Value = "'-1.48%"
$ws.Range("D18").Value = "'2.131"
$ws.Range("E18").Value = "'3.56%"
$ws.Range("D19").Value = "'0.3378"
$ws.Range("E19").Value = "'1.36%"
$ws.Range("D20").Value = "'0.1326"
$ws.Range("D21").Value = "'4.330"
$ws.Range("E21").Value = "'3.88%"
$ws.Range("D22").Value = "'0.2225"
$ws.Range("E22").Value = "'5.85%"
$ws.Range("D23").Value = "'0.04569"
$ws.Range("E23").Value = "'0.60%"
$ws.Range("D24").Value = "'0.001217"
$ws.Range("E24").Value = "'0.64%"
$ws.Range("D25").Value = "'0.004442"
$ws.Range("E25").Value = "'6.82%"
$ws.Range("E26").Value = "'-0.24%"
$ws.Range("E27").Value = "'8.75%"
$ws.Range("D39").Value = "'0.01779"
$ws.Range("E39").Value = "'13.43%"
$ws.Range("D40").Value = "'0.04798"
$ws.Range("E40").Value = "'8.10%"
$ws.Range("D41").Value = "'0.007731"
$ws.Range("E41").Value = "'4.95%"
$ws.Range("E42").Value = "'6.12%"
$ws.Range("D43").Value = "'0.007125"
$ws.Range("E43").Value = "'-20.81%"
$ws.Range("E44").Value = "'-2.08%"
$ws.Range("D45").Value = "'0.009552"
$ws.Range("E45").Value = "'4.85%"
$ws.Range("D46").Value = "'0.00005923"
$ws.Range("E46").Value = "'-3.17%"
$ws.Range("E47").Value = "'-0.23%"
$ws.Range("E48").Value = "'26.36%"
$ws.Range("E49").Value = "'-0.12%"
$ws.Range("D50").Value = "'0.00002097"
$ws.Range("E50").Value = "'-0.23%"
$ws.Range("D51").Value = "'0.0001997"
$ws.Range("E51").Value = "'-0.23%"
